# Append a new login-info row (esti / esti19 / ighfhgd / e@123456) to the
# "Table1" worksheet, right after the existing last row (row 4), mirroring
# the other rows already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "esti"
$ws.Range("B5").Value = "esti19"
$ws.Range("C5").Value = "ighfhgd"
$ws.Range("D5").Value = "e@123456"
